$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-parsed as a number
# by Excel (losing the literal "x.xx"/"x.xx.xx" display, e.g. "1.00" -> 1,
# "0.0681" -> 6.81E-2). Force them to Text first so the literal string is
# preserved exactly, matching the original inline-string cell content.
$textForceCells = @("D4", "D5", "D6", "D14", "D15", "D19", "D21", "D24", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D43", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
$ws.Range("D2").Value = '57.164.66'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '3.014.38'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '516.83'
$ws.Range("E5").Value = '  +4.34%  '
$ws.Range("D6").Value = '139.87'
$ws.Range("E6").Value = '  +4.98%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +3.38%  '
$ws.Range("E9").Value = '  +5.25%  '
$ws.Range("E10").Value = '  +6.19%  '
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").Value = '3.529.81'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").Value = '25.92'
$ws.Range("E14").Value = '  +4.45%  '
$ws.Range("D15").Value = '0.0000160'
$ws.Range("E15").Value = '  +11.48%  '
$ws.Range("D16").Value = '57.152.23'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = '3.019.23'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '12.67'
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("E20").Value = '  +2.82%  '
$ws.Range("D21").Value = '329.89'
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  +5.20%  '
$ws.Range("D24").Value = '63.83'
$ws.Range("E24").Value = '  +4.96%  '
$ws.Range("E25").Value = '  +5.27%  '
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  +5.89%  '
$ws.Range("E28").Value = '  +2.80%  '
$ws.Range("D29").Value = '7.18'
$ws.Range("E29").Value = '  +7.81%  '
$ws.Range("E30").Value = '  +5.90%  '
$ws.Range("D31").Value = '1.22'
$ws.Range("E31").Value = '  +4.52%  '
$ws.Range("D32").Value = '20.73'
$ws.Range("E32").Value = '  +5.11%  '
$ws.Range("D33").Value = '157.89'
$ws.Range("E33").Value = '  +4.59%  '
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").Value = '  +4.39%  '
$ws.Range("D35").Value = '5.75'
$ws.Range("E35").Value = '  +1.60%  '
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '0.0681'
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("D38").Value = '24.26'
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("D39").Value = '3.045.94'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("D40").Value = '37.33'
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = '2.296.87'
$ws.Range("E42").Value = '  +5.99%  '
$ws.Range("D43").Value = '0.650'
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("E44").Value = '  +5.16%  '
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("D48").Value = '0.0241'
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("D49").Value = '5.88'
$ws.Range("E49").Value = '  +5.78%  '
$ws.Range("D50").Value = '19.39'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").Value = '0.0881'
$ws.Range("E51").Value = '  +4.01%  '

# Restore the default cell style: the Text number-format coercion above
# must not leave a visible style change on the cells themselves.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
